$d = $word.ActiveDocument

# 1) Title paragraph: update the date and add a new line with the paper title.
$d.Paragraphs(1).Range.Text = "המאמר היומי של מייק: 19.07.25`vGENARM: Reward Guided Generation with Autoregressive Reward Model for Test-Time Alignment"

# 2) Second paragraph: replace "Reinforcement Pre-Training" with the new intro text.
$d.Paragraphs(2).Range.Text = "עבר כבר שבוע מהסקירה האחרונה והרגשתי צורך דחוף לסקור איזה מאמר. האמת די הרבה זמן לא היתה לי הפסקה כזו גדולה לצערי גם רוחב הפס שלי אינו אינסופי. טוב, יאללה מתחילים לסקור."

# 3) Third paragraph.
$d.Paragraphs(3).Range.Text = "המאמר מדבר על גנרוט דאטה באמצעות מודל שפה תוך התחשבות(כוונון) במודל reward (תגמול) חיצוני האומד את איכות הטקסט המגונרט. האיכות נאמדת רק כאשר הגנרוט נגמר בסוף הטקסט כלומר עבור התשובה המלאה. נציין כי ניתן להשתמש בטריק שידוע לנו משיטת DPO שזה Direct Preference Optimization כדי להתחשב בציון ממודל התגמול עבור התשובה כדי להכווין את התפלגות הגנרוט של המודל. "

# 4) Fourth paragraph.
$d.Paragraphs(4).Range.Text = "שיטת תיקון ממודל התגמול נובעת מהנוסחה עבור פונקציית לוס של אימון מודל שפה עם LHF כאשר המטרה (של האימון) היא למקסם את הציונים של התשובות של המודל עם רגולריזציה שמנסה לשמור את התפלגות המודל המאומן קרובה להתפלגות ההתחלתית של המודל במונחי מרחק KL. בד""כ מבצעים אימון כזה על דאטהסט של שאלות עם תשובות רצויות ולא רצויות שכאמור למקסם את יחס הציונים ביניהם. התיקון מתבצע ללוג של הסתברות של משוב המגונרט המלא y (בהינתן ההקשר x) על ידי החיבור של הציון (r(x, y (ממושקל) ופונקציית נרמול התלויה ב x בלבד (המאמר לא מרחיב על אופן שערכו). "

# 5) Fifth paragraph.
$d.Paragraphs(5).Range.Text = "אז איך כל הסיפור הזה (גנרוט עם פונקציית תגמול ללא אימון RLHF) עבד לפני המאמר הזה? בזמן הגנרוט בשביל לגנרט טוקו הבא בהינתן הטוקנים שכבר גונרטו אנו דוגמים כמה המשכים עד סוף התשובה ואז ניתן להשתמש במודל תגמול בשביל לשערך את איכותה. אז הטוקן שהוא נמצא ההשלמה בעלת הנראות המתוקנת הגבוהה ביותר. הסיבה לכך היא העובדה שלא ניתן לשערך רק את התשובה המלאה ולא חלקית שלא מאפשרת חישוב התיקון עבור כל טוקו מגנט בצורה ישירה. יש עוד שיטות לעשות את זה אבל הם או לא יעילות או ביצועיהן לא כאלו טובות."

# 6) Sixth paragraph.
$d.Paragraphs(6).Range.Text = "המאמר המסוקר מציע לאמן מודל שמטרתו היא לשערך (r(x, y עבור תשובות חלקיות בהתבסס על הדאטהסט של שאלות עם תשובות רצויות ולא רצויות. המאמר ממקסם את היחס סכום התגמולים עבור כל הטוקנים של התשובות הרצויות לאלו של לא רצויות. מודל זה כמובן מתבסס על מודל שפה עם ראש מאומן כמו שנעשה עבור אימון מודל תגמול רגיל עבור תשובות מלאות. המחברים טוענים כי מודל תגמול קטן יחסית למשל 7B מסוגל לשפר את איכות הגנרוט בהתאם ל alignment הרצוי עבור גנרוט למודל הרבה יותר גדול עם 70B פרמטרים. "

# 7) Seventh paragraph.
$d.Paragraphs(7).Range.Text = "בצורה כזו ניתן לבצע גנרוט בהתבסס על כמה מדיניות הalignment שכל אחת מיוצגת על ידי דאטהסט משלה. אחרי אימון של מודל התגמול עבור כל אחת מהן ניתן לבנות את התיקון ללוג של נראות עבור טוקו הבא על ידי סכום ממושקל של התגמולים עבור כל אחת מהן כאשר המשקול תלוי במידת התחשבות בכל אחת מהמדינות alignment אלו."

# 8) Eighth paragraph (was the arxiv link) becomes a short closing remark.
$d.Paragraphs(8).Range.Text = "מאמר קליל אך עם זאת די מעניין"

# 9) New final paragraph with the updated arxiv link.
$d.Paragraphs(8).Range.InsertParagraphAfter()
$d.Paragraphs(9).Range.Text = "  https://arxiv.org/abs/2410.0819"
